# Insert a new data row at row 30 (pushing the existing rows 30-102 down to
# 31-103) and populate it with a new price observation.
# This grows the used range from A1:R102 to A1:R103.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(30).Insert()

$ws.Range("A30").Value = 4
$ws.Range("B30").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C30").Value = "Los Lagos"
$ws.Range("D30").Value = "2022-02-28"
$ws.Range("E30").Value = 10
$ws.Range("F30").Value = 100112052
$ws.Range("G30").Value = "Albahaca"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 70
$ws.Range("K30").Value = 6000
$ws.Range("L30").Value = 6000
$ws.Range("M30").Value = 6000
$ws.Range("N30").Value = '$/docena de matas'
$ws.Range("O30").Value = "Región Metropolitana"
$ws.Range("P30").Value = 1000
$ws.Range("Q30").Value = 6
$ws.Range("R30").Value = "Hortaliza"
